# Cleanup cases: add a new "ug" column to the EDTSlot sheet, with value
# "1,1,1,1" for every data row, and leave the workbook with EDTSlot as the
# active/selected sheet (matching the author's final on-screen state).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. EDTSlot sheet: add column E ("ug") filled with "1,1,1,1"
# ---------------------------------------------------------------------
$edt = $wb.Worksheets.Item("EDTSlot")

$edt.Range("E1").Value = "ug"
for ($row = 2; $row -le 25; $row++) {
    $edt.Cells.Item($row, 5).Value = "1,1,1,1"
}

# ---------------------------------------------------------------------
# 2. Make EDTSlot the active sheet / selection, moving it off GCost
#    (GCost was previously the selected tab).
# ---------------------------------------------------------------------
$edt.Activate()
$edt.Range("N17").Select()
